$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (existing rows 52.. shift down to 53..)
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with its two visible values
$ws.Range("A52").Value = "tbl_summary"
$ws.Range("B52").Value = "range"

# The source row only carried A:F formatting (C:F stayed blank/untouched);
# clear G:M so those cells don't linger with inherited formatting/content
$ws.Range("G52:M52").Clear()

Write-Output "done"
